# AddressMasterTemplate.xlsx - "some changes regarding to container-master-upload"
#
# The header label in column B was renamed from "addressCode" to "address",
# and the active selection was left on B1 (the cell that was just edited)
# instead of the old F11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "addressCode" header to "address"
$ws.Range("B1").Value = "address"

# Leave the selection on the cell that was edited
$ws.Range("B1").Select()
